$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 56
$ws.Cells.Item($row, 1).Value = 45702
$ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat
$ws.Cells.Item($row, 2).Value = "continuing with the area"
$ws.Cells.Item($row, 3).Value = 4

$null = $ws.Range("C58").Select()
